$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 10

$ws.Cells.Item($row, 1).Value = 111906849
$ws.Cells.Item($row, 2).Value = 88967
$ws.Cells.Item($row, 3).Value = "Ovaliderad"
$ws.Cells.Item($row, 4).Value = "DD"
$ws.Cells.Item($row, 5).Value = 6039940
$ws.Cells.Item($row, 6).Value = "Mandarinfingersvamp"
$ws.Cells.Item($row, 7).Value = "Ramaria tridentina"
$ws.Cells.Item($row, 8).Value = "Schild"
$ws.Cells.Item($row, 9).Value = ""
$ws.Cells.Item($row, 11).Value = ""
$ws.Cells.Item($row, 16).Value = "Torrkölen (Torrkölen), Nb"
$ws.Cells.Item($row, 17).Value = 813178.8074009671
$ws.Cells.Item($row, 18).Value = 7316199.822832054
$ws.Cells.Item($row, 19).Value = 20
$ws.Cells.Item($row, 20).Value = "Norrbotten"
$ws.Cells.Item($row, 21).Value = "Boden"
$ws.Cells.Item($row, 22).Value = "Norrbotten"
$ws.Cells.Item($row, 23).Value = "Överluleå"
$ws.Range("Y" + $row + ":AB" + $row).NumberFormat = "@"
$ws.Cells.Item($row, 25).Value = "2023-09-05"
$ws.Cells.Item($row, 26).Value = "11:11"
$ws.Cells.Item($row, 27).Value = "2023-09-05"
$ws.Cells.Item($row, 28).Value = "11:11"
$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false
$ws.Cells.Item($row, 33).Value = $false
$ws.Cells.Item($row, 46).Value = ""
$ws.Cells.Item($row, 49).Value = "Linnea Åsedahl"
$ws.Cells.Item($row, 50).Value = "Linnea Åsedahl"
$ws.Cells.Item($row, 51).Value = ""
